$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state for rows 2-46: column B (id) is derived from column C (speaker_variant)
# as "#" + lowercase(C) with spaces replaced by hyphens. Column D (is_prefered) is
# cleared for every row. Rows are also reordered per the source diff.

$ws.Cells.Item(2, 2).Value = '#donnare-turbanta?-donnare-turbanta?'
$ws.Cells.Item(2, 3).Value = 'Donnare Turbanta? Donnare Turbanta?'
$ws.Cells.Item(2, 4).Value = ''

$ws.Cells.Item(3, 2).Value = '#scherm'
$ws.Cells.Item(3, 3).Value = 'Scherm'
$ws.Cells.Item(3, 4).Value = ''

$ws.Cells.Item(4, 2).Value = '#marg'
$ws.Cells.Item(4, 3).Value = 'Marg'
$ws.Cells.Item(4, 4).Value = ''

$ws.Cells.Item(5, 2).Value = '#luc'
$ws.Cells.Item(5, 3).Value = 'Luc'
$ws.Cells.Item(5, 4).Value = ''

$ws.Cells.Item(6, 2).Value = '#leerl'
$ws.Cells.Item(6, 3).Value = 'Leerl'
$ws.Cells.Item(6, 4).Value = ''

$ws.Cells.Item(7, 2).Value = '#nosta-panna,-stara-ba,-ba?'
$ws.Cells.Item(7, 3).Value = 'Nosta panna, Stara ba, ba?'
$ws.Cells.Item(7, 4).Value = ''

$ws.Cells.Item(8, 2).Value = '#zangm.'
$ws.Cells.Item(8, 3).Value = 'Zangm.'
$ws.Cells.Item(8, 4).Value = ''

$ws.Cells.Item(9, 2).Value = '#tinoster-turba'
$ws.Cells.Item(9, 3).Value = 'Tinoster Turba'
$ws.Cells.Item(9, 4).Value = ''

$ws.Cells.Item(10, 2).Value = '#donnare-turbanta,-donnare-turbanta'
$ws.Cells.Item(10, 3).Value = 'Donnare Turbanta, Donnare Turbanta'
$ws.Cells.Item(10, 4).Value = ''

$ws.Cells.Item(11, 2).Value = '#knecht'
$ws.Cells.Item(11, 3).Value = 'Knecht'
$ws.Cells.Item(11, 4).Value = ''

$ws.Cells.Item(12, 2).Value = '#mahometa-par-jourdina,-mipregersera-mesina'
$ws.Cells.Item(12, 3).Value = 'Mahometa par Jourdina, Mipregersera Mesina'
$ws.Cells.Item(12, 4).Value = ''

$ws.Cells.Item(13, 2).Value = '#lucil'
$ws.Cells.Item(13, 3).Value = 'Lucil'
$ws.Cells.Item(13, 4).Value = ''

$ws.Cells.Item(14, 2).Value = '#pag'
$ws.Cells.Item(14, 3).Value = 'Pag'
$ws.Cells.Item(14, 4).Value = ''

$ws.Cells.Item(15, 2).Value = '#tistar-nobile.-non-ster-ba-bo-la,-piglisa-bo-la'
$ws.Cells.Item(15, 3).Value = 'Tistar Nobile. Non ster ba bo la, Piglisa bo la'
$ws.Cells.Item(15, 4).Value = ''

$ws.Cells.Item(16, 2).Value = '#iour'
$ws.Cells.Item(16, 3).Value = 'Iour'
$ws.Cells.Item(16, 4).Value = ''

$ws.Cells.Item(17, 2).Value = '#cleo'
$ws.Cells.Item(17, 3).Value = 'Cleo'
$ws.Cells.Item(17, 4).Value = ''

$ws.Cells.Item(18, 2).Value = '#dorim'
$ws.Cells.Item(18, 3).Value = 'Dorim'
$ws.Cells.Item(18, 4).Value = ''

$ws.Cells.Item(19, 2).Value = '#zangm'
$ws.Cells.Item(19, 3).Value = 'Zangm'
$ws.Cells.Item(19, 4).Value = ''

$ws.Cells.Item(20, 2).Value = '#lurd'
$ws.Cells.Item(20, 3).Value = 'Lurd'
$ws.Cells.Item(20, 4).Value = ''

$ws.Cells.Item(21, 2).Value = '#dor'
$ws.Cells.Item(21, 3).Value = 'Dor'
$ws.Cells.Item(21, 4).Value = ''

$ws.Cells.Item(22, 2).Value = '#philos'
$ws.Cells.Item(22, 3).Value = 'Philos'
$ws.Cells.Item(22, 4).Value = ''

$ws.Cells.Item(23, 2).Value = '#graaf'
$ws.Cells.Item(23, 3).Value = 'Graaf'
$ws.Cells.Item(23, 4).Value = ''

$ws.Cells.Item(24, 2).Value = '#iourd'
$ws.Cells.Item(24, 3).Value = 'Iourd'
$ws.Cells.Item(24, 4).Value = ''

$ws.Cells.Item(25, 2).Value = '#schermm'
$ws.Cells.Item(25, 3).Value = 'Schermm'
$ws.Cells.Item(25, 4).Value = ''

$ws.Cells.Item(26, 2).Value = '#non-star-turbanta'
$ws.Cells.Item(26, 3).Value = 'Non star Turbanta'
$ws.Cells.Item(26, 4).Value = ''

$ws.Cells.Item(27, 2).Value = '#chira,-ba,-ba.-chira,-ba,-ba'
$ws.Cells.Item(27, 3).Value = 'Chira, ba, ba. Chira, ba, ba'
$ws.Cells.Item(27, 4).Value = ''

$ws.Cells.Item(28, 2).Value = '#dara,-dara,-bastonara,-bastonara'
$ws.Cells.Item(28, 3).Value = 'Dara, dara, bastonara, bastonara'
$ws.Cells.Item(28, 4).Value = ''

$ws.Cells.Item(29, 2).Value = '#jourd.'
$ws.Cells.Item(29, 3).Value = 'Jourd.'
$ws.Cells.Item(29, 4).Value = ''

$ws.Cells.Item(30, 2).Value = '#dansm'
$ws.Cells.Item(30, 3).Value = 'Dansm'
$ws.Cells.Item(30, 4).Value = ''

$ws.Cells.Item(31, 2).Value = '#zang'
$ws.Cells.Item(31, 3).Value = 'Zang'
$ws.Cells.Item(31, 4).Value = ''

$ws.Cells.Item(32, 2).Value = '#non-tenar-bonta-questa,-star-ultima-fronta'
$ws.Cells.Item(32, 3).Value = 'Non tenar bonta questa, star ultima Fronta'
$ws.Cells.Item(32, 4).Value = ''

$ws.Cells.Item(33, 2).Value = '#phil'
$ws.Cells.Item(33, 3).Value = 'Phil'
$ws.Cells.Item(33, 4).Value = ''

$ws.Cells.Item(34, 2).Value = '#nicol'
$ws.Cells.Item(34, 3).Value = 'Nicol'
$ws.Cells.Item(34, 4).Value = ''

$ws.Cells.Item(35, 2).Value = '#se-tier-sabier,-tires-pondier,-senonsahir'
$ws.Cells.Item(35, 3).Value = 'Se tier Sabier, Tires pondier, senonsahir'
$ws.Cells.Item(35, 4).Value = ''

$ws.Cells.Item(36, 2).Value = '#starrabon-turca-jourdina?'
$ws.Cells.Item(36, 3).Value = 'Starrabon Turca Jourdina?'
$ws.Cells.Item(36, 4).Value = ''

$ws.Cells.Item(37, 2).Value = '#kn'
$ws.Cells.Item(37, 3).Value = 'Kn'
$ws.Cells.Item(37, 4).Value = ''

$ws.Cells.Item(38, 2).Value = '#margo'
$ws.Cells.Item(38, 3).Value = 'Margo'
$ws.Cells.Item(38, 4).Value = ''

$ws.Cells.Item(39, 2).Value = '#ivolla,-ivolla,-ivolla'
$ws.Cells.Item(39, 3).Value = 'Ivolla, ivolla, ivolla'
$ws.Cells.Item(39, 4).Value = ''

$ws.Cells.Item(40, 2).Value = '#cov'
$ws.Cells.Item(40, 3).Value = 'Cov'
$ws.Cells.Item(40, 4).Value = ''

$ws.Cells.Item(41, 2).Value = '#jourd'
$ws.Cells.Item(41, 3).Value = 'Jourd'
$ws.Cells.Item(41, 4).Value = ''

$ws.Cells.Item(42, 2).Value = '#snij'
$ws.Cells.Item(42, 3).Value = 'Snij'
$ws.Cells.Item(42, 4).Value = ''

$ws.Cells.Item(43, 2).Value = '#phlos'
$ws.Cells.Item(43, 3).Value = 'Phlos'
$ws.Cells.Item(43, 4).Value = ''

$ws.Cells.Item(44, 2).Value = '#nic'
$ws.Cells.Item(44, 3).Value = 'Nic'
$ws.Cells.Item(44, 4).Value = ''

$ws.Cells.Item(45, 2).Value = '#io'
$ws.Cells.Item(45, 3).Value = 'Io'
$ws.Cells.Item(45, 4).Value = ''

$ws.Cells.Item(46, 2).Value = '#cleont'
$ws.Cells.Item(46, 3).Value = 'Cleont'
$ws.Cells.Item(46, 4).Value = ''
